$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coal heating technology (id_heating_technology = 45) availability turned off
# (1 -> 0) for years 2015-2050 (columns K:AT) across the three heating-system-action rows.
$ws.Range("K24:AT24").Value = 0
$ws.Range("K49:AT49").Value = 0
$ws.Range("K74:AT74").Value = 0

# Update the active selection left by the author before saving.
$ws.Range("U4").Select()
